# Append a new scraped entry (2025-08-27 18:25 JST) to the top of the
# "ランサーズ" listing sheet, and a matching stats row to "統計".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ランサーズ
$ws2 = $wb.Worksheets.Item(2)   # 統計

# ---------------------------------------------------------------------
# Sheet 1 (ランサーズ): insert a brand-new row 2, pushing the existing
# 11 data rows (old rows 2-12) down to rows 3-13. Rows.Insert() shifts
# cell values/styles correctly, but it does NOT move the worksheet's
# Hyperlinks collection, so those are rebuilt from scratch afterward.
# ---------------------------------------------------------------------
$ws1.Rows.Item(2).Insert()

$ws1.Cells.Item(2,1).Value = "2025-08-27 18:25:32"
$ws1.Cells.Item(2,2).Value = "Qt / C++ ベースのWindowsアプリとデバイスファームウェアの修正"
$ws1.Cells.Item(2,3).Value = "システム開発"
$ws1.Cells.Item(2,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws1.Cells.Item(2,5).Value = "期限情報なし"
$ws1.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5380896"
$ws1.Cells.Item(2,7).Value = 45
$ws1.Cells.Item(2,8).Value = "◇アプリ"

# Rebuild the hyperlinks for column F (rows 2-13) in document order so
# the relationship ids end up rId1..rId12, in the same top-to-bottom
# order the rows appear.
$ws1.Hyperlinks.Delete()

$urls = @(
    "https://www.lancers.jp/work/detail/5380896",
    "https://www.lancers.jp/work/detail/5380830",
    "https://www.lancers.jp/work/detail/5016989",
    "https://www.lancers.jp/work/detail/5273634",
    "https://www.lancers.jp/work/detail/5217096",
    "https://www.lancers.jp/work/detail/5380343",
    "https://www.lancers.jp/work/detail/5380337",
    "https://www.lancers.jp/work/detail/5380683",
    "https://www.lancers.jp/work/detail/5380747",
    "https://www.lancers.jp/work/detail/5341051",
    "https://www.lancers.jp/work/detail/5380357",
    "https://www.lancers.jp/work/detail/5380420"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $cell = $ws1.Cells.Item($row, 6)
    $ws1.Hyperlinks.Add($cell, $urls[$i])
    $cell.Style = "Hyperlink"
}

# Column D got a bit wider in this revision (raw OOXML width 28 -> 30).
# Excel's ColumnWidth property is offset from the raw stored width by
# 5/6 of a character, so back that out to land exactly on 30.
$ws1.Columns.Item(4).ColumnWidth = 30 - 5/6

# ---------------------------------------------------------------------
# Sheet 2 (統計): append the matching stats row (old last row + 1).
# ---------------------------------------------------------------------
$ws2.Cells.Item(5,1).Value = "2025-08-27T18:25:32.305687"
$ws2.Cells.Item(5,2).Value = 12
$ws2.Cells.Item(5,3).Value = "全案件リスト"
$ws2.Cells.Item(5,4).Value = 75
$ws2.Cells.Item(5,5).Value = 3
$ws2.Cells.Item(5,6).Value = 6
$ws2.Cells.Item(5,7).Value = 12
